# Insert a new data row before existing row 70 (shifts rows 70-133 down to 71-134)
# and populate it with a new weekly price observation (dimension grows to A1:R134).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(70).Insert()

$ws.Range("A70").Value = 8
$ws.Range("B70").Value = "Terminal La Palmera de La Serena"
$ws.Range("C70").Value = "Coquimbo"
$ws.Range("D70").Value = 44651
$ws.Range("E70").Value = 4
$ws.Range("F70").Value = 100112044
$ws.Range("G70").Value = "Perejil"
$ws.Range("H70").Value = "Sin especificar"
$ws.Range("I70").Value = "Primera"
$ws.Range("J70").Value = 2500
$ws.Range("K70").Value = 2000
$ws.Range("L70").Value = 2500
$ws.Range("M70").Value = 2250
$ws.Range("N70").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O70").Value = "Provincia del Elquí"
$ws.Range("P70").Value = 1500
$ws.Range("Q70").Value = 1.5
$ws.Range("R70").Value = "Hortaliza"
